# Automatic update of files.
# - Bumps the "Förändrad" (C) column date for every data row from 2023-09-15
#   (45184) to 2023-09-17 (45186).
# - Adds a friendly display-text second argument to the HYPERLINK() formulas
#   in columns S, T, V, W, X, Y (only present for rows 2-12), using the
#   "Beteckning" value from column A of the same row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "Förändrad" date column (C2:C45) -----------------------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$ws.Range("C2:C" + $lastRow).Value2 = 45186

# --- 2. Add display text to HYPERLINK formulas (columns S, T, V, W, X, Y) -
#    Only rows that still have a plain single-argument HYPERLINK("url")
#    formula are touched; rows without links (or already-updated links) are
#    left alone.
$hyperlinkCols = @("S", "T", "V", "W", "X", "Y")

for ($row = 2; $row -le $lastRow; $row++) {
    $beteckning = $ws.Range("A" + $row).Value2

    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Range($col + $row)
        $formula = $cell.Formula
        if ([string]::IsNullOrEmpty($formula)) {
            continue
        }

        # Only touch cells that are still plain single-argument HYPERLINK
        # formulas, i.e. they don't already carry a display-text argument.
        if ($formula -match '^=HYPERLINK\("([^"]*)"\)$') {
            $url = $matches[1]
            $cell.Formula = '=HYPERLINK("' + $url + '", "' + $beteckning + '")'
        }
    }
}

Write-Output "done"
